$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the room ("Salle") column F for each course row, which was
# previously left blank.
$ws.Range("F3").Value  = "U3-110"
$ws.Range("F6").Value  = "U3-4"
$ws.Range("F9").Value  = "U3-4"
$ws.Range("F12").Value = "U3-Amphi"
$ws.Range("F14").Value = "U3-4"
$ws.Range("F16").Value = "U3-Amphi"
$ws.Range("F19").Value = "U3-Amphi"
$ws.Range("F20").Value = "U3-Amphi"
$ws.Range("F23").Value = "U3-Amphi"
